$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B1 header text (was "Cell_type") to the new instructional text
$ws.Range("B1").Value = "Add annotations used in your experiment per column"

# Remove the old "Condition" / "Mouse" / "Passage" columns (C1:E1) entirely
$ws.Range("C1:E1").ClearContents()

# Select B2 to match the workbook's saved selection state
$ws.Range("B2").Select()
